$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 219, shifting existing rows 219:297 down to 220:298,
# then populate the new row with the new weekly price entry.
$ws.Rows("219:219").Insert()

$ws.Range("A219").Value2 = 10
$ws.Range("B219").Value2 = "Vega Modelo de Temuco"
$ws.Range("C219").Value2 = "La Araucanía"
$ws.Range("D219").Value2 = 45141
$ws.Range("E219").Value2 = 9
$ws.Range("F219").Value2 = 100112012
$ws.Range("G219").Value2 = "Espinaca"
$ws.Range("H219").Value2 = "Sin especificar"
$ws.Range("I219").Value2 = "Primera"
$ws.Range("J219").Value2 = 120
$ws.Range("K219").Value2 = 8000
$ws.Range("L219").Value2 = 8000
$ws.Range("M219").Value2 = 8000
$ws.Range("N219").Value2 = "$/docena de paquetes"
$ws.Range("O219").Value2 = "Región de La Araucanía"
$ws.Range("P219").Value2 = 667
$ws.Range("Q219").Value2 = 12
$ws.Range("R219").Value2 = "Hortaliza"

# Match the date cell style used by the rest of column D (datetime number format).
$ws.Range("D219").NumberFormat = $ws.Range("D220").NumberFormat
